$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 97; $r++) {
    $ws.Cells.Item($r, 1).Value = "q" + ($r - 2)
}
